# Generate Report for Handoff
# Status moves from "In Translation" to "Ready for handoff", and the
# handoff timestamps advance to reflect the newly generated report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Sheets.Item("Overview")
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-03-25 12:44:55"

$zhcn = $wb.Sheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-25 12:44:50"

$dede = $wb.Sheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-25 12:44:55"
